# Auto-upload VRF Excel file
# Appends two new data rows (row 2 and row 3) to the "new" worksheet,
# which previously contained only the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("new")

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = "wefsdf"
$ws.Cells.Item(2, 2).Value = 1
# Force "234234" to be stored as text (not a number), same as typing '234234 in Excel.
$ws.Cells.Item(2, 3).Value = "'234234"
# D2, E2, F2 are left blank (empty cells) for this row.

# --- Row 3 ---
# A3, B3, C3 are left blank (empty cells) for this row.
$ws.Cells.Item(3, 4).Value = "wewdrewsfdesw"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = "sdfcsdsdfssdf, sfwerwr"
